$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    $origStyle = $Cell.Style
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Range("D2") "71.058.89"
Set-TextValue $ws.Range("E2") "  +4.65%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.621.92"
Set-TextValue $ws.Range("E3") "  +5.41%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.08%  "

# Row 5
Set-TextValue $ws.Range("D5") "606.03"
Set-TextValue $ws.Range("E5") "  +2.98%  "

# Row 6
Set-TextValue $ws.Range("D6") "180.93"
Set-TextValue $ws.Range("E6") "  +3.68%  "

# Row 7
Set-TextValue $ws.Range("E7") "  +0.02%  "

# Row 8
Set-TextValue $ws.Range("E8") "  +1.94%  "

# Row 9
Set-TextValue $ws.Range("D9") "2.621.94"
Set-TextValue $ws.Range("E9") "  +5.40%  "

# Row 10
Set-TextValue $ws.Range("E10") "  +14.60%  "

# Row 12
Set-TextValue $ws.Range("E12") "  +3.37%  "

# Row 13
Set-TextValue $ws.Range("D13") "5.04"
Set-TextValue $ws.Range("E13") "  +1.53%  "

# Row 14
Set-TextValue $ws.Range("D14") "3.097.81"
Set-TextValue $ws.Range("E14") "  +5.35%  "

# Row 15
Set-TextValue $ws.Range("D15") "26.67"
Set-TextValue $ws.Range("E15") "  +5.29%  "

# Row 16
Set-TextValue $ws.Range("E16") "  +7.39%  "

# Row 17
Set-TextValue $ws.Range("D17") "71.036.63"
Set-TextValue $ws.Range("E17") "  +4.77%  "

# Row 18
Set-TextValue $ws.Range("D18") "2.617.69"
Set-TextValue $ws.Range("E18") "  +5.92%  "

# Row 19
Set-TextValue $ws.Range("D19") "380.48"
Set-TextValue $ws.Range("E19") "  +9.71%  "

# Row 20
Set-TextValue $ws.Range("E20") "  +7.09%  "

# Row 21
Set-TextValue $ws.Range("E21") "  +6.38%  "

# Row 22
Set-TextValue $ws.Range("E22") "  +1.18%  "

# Row 23
Set-TextValue $ws.Range("D23") "72.01"
Set-TextValue $ws.Range("E23") "  +1.75%  "

# Row 24
Set-TextValue $ws.Range("E24") "  +6.47%  "

# Row 25
Set-TextValue $ws.Range("E25") "  -0.04%  "

# Row 26
Set-TextValue $ws.Range("D26") "1.85"
Set-TextValue $ws.Range("E26") "  +9.79%  "

# Row 27
Set-TextValue $ws.Range("D27") "9.65"
Set-TextValue $ws.Range("E27") "  +9.15%  "

# Row 29
Set-TextValue $ws.Range("E29") "  -0.01%  "

# Row 30
Set-TextValue $ws.Range("D30") "0.0₃0952"
Set-TextValue $ws.Range("E30") "  +7.07%  "

# Row 31
Set-TextValue $ws.Range("D31") "528.44"
Set-TextValue $ws.Range("E31") "  +5.88%  "

# Row 32
Set-TextValue $ws.Range("D32") "8.01"
Set-TextValue $ws.Range("E32") "  +3.76%  "

# Row 33
Set-TextValue $ws.Range("D33") "1.33"
Set-TextValue $ws.Range("E33") "  +6.61%  "

# Row 34
Set-TextValue $ws.Range("E34") "  +4.07%  "

# Row 35
Set-TextValue $ws.Range("D35") "1.00"
Set-TextValue $ws.Range("E35") "  +0.04%  "

# Row 36
Set-TextValue $ws.Range("D36") "164.03"
Set-TextValue $ws.Range("E36") "  -0.34%  "

# Row 37
Set-TextValue $ws.Range("E37") "  -0.66%  "

# Row 38
Set-TextValue $ws.Range("D38") "19.12"
Set-TextValue $ws.Range("E38") "  +4.63%  "

# Row 39
Set-TextValue $ws.Range("B39") "Stacks"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D39") "1.87"
Set-TextValue $ws.Range("E39") "  +8.34%  "

# Row 40
Set-TextValue $ws.Range("B40") "WhiteBITCoin"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws.Range("D40") "18.94"
Set-TextValue $ws.Range("E40") "  +1.65%  "

# Row 41
Set-TextValue $ws.Range("D41") "1.38"
Set-TextValue $ws.Range("E41") "  +5.38%  "

# Row 42
Set-TextValue $ws.Range("E42") "  +0.09%  "

# Row 43
Set-TextValue $ws.Range("B43") "dogwifhat"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D43") "2.61"
Set-TextValue $ws.Range("E43") "  +9.38%  "

# Row 44
Set-TextValue $ws.Range("B44") "RenderToken"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue $ws.Range("D44") "5.02"
Set-TextValue $ws.Range("E44") "  +5.25%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.333"
Set-TextValue $ws.Range("E45") "  +2.52%  "

# Row 46
Set-TextValue $ws.Range("D46") "40.13"
Set-TextValue $ws.Range("E46") "  +3.86%  "

# Row 47
Set-TextValue $ws.Range("D47") "153.83"
Set-TextValue $ws.Range("E47") "  +4.17%  "

# Row 48
Set-TextValue $ws.Range("B48") "Filecoin"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D48") "3.65"
Set-TextValue $ws.Range("E48") "  +3.71%  "

# Row 49
Set-TextValue $ws.Range("B49") "BabyDogeCoin"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D49") "0.0₆0273"
Set-TextValue $ws.Range("E49") "  +7.99%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.533"
Set-TextValue $ws.Range("E50") "  +4.20%  "

# Row 51
Set-TextValue $ws.Range("D51") "1.67"
Set-TextValue $ws.Range("E51") "  +7.10%  "
